$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "26.065.57"
Set-TextValue $ws.Range("E2") "  -1.96%  "
Set-TextValue $ws.Range("D3") "1.668.32"
Set-TextValue $ws.Range("E3") "  -1.59%  "
Set-TextValue $ws.Range("E4") "  -0.07%  "
Set-TextValue $ws.Range("D5") "216.94"
Set-TextValue $ws.Range("E5") "  -0.95%  "
Set-TextValue $ws.Range("D6") "0.5109"
Set-TextValue $ws.Range("E6") "  +0.52%  "
Set-TextValue $ws.Range("D7") "1.005"
Set-TextValue $ws.Range("E7") "  -0.09%  "
Set-TextValue $ws.Range("D8") "0.2655"
Set-TextValue $ws.Range("E8") "  +0.33%  "
Set-TextValue $ws.Range("D9") "0.06413"
Set-TextValue $ws.Range("E9") "  +2.04%  "
Set-TextValue $ws.Range("D10") "21.91"
Set-TextValue $ws.Range("E10") "  -0.86%  "
Set-TextValue $ws.Range("D11") "0.07431"
Set-TextValue $ws.Range("E11") "  +0.88%  "
Set-TextValue $ws.Range("D12") "1.670.71"
Set-TextValue $ws.Range("E12") "  -1.47%  "
Set-TextValue $ws.Range("D13") "4.499"
Set-TextValue $ws.Range("E13") "  -0.23%  "
Set-TextValue $ws.Range("D14") "0.5864"
Set-TextValue $ws.Range("E14") "  +0.54%  "
Set-TextValue $ws.Range("D15") "0.000008568"
Set-TextValue $ws.Range("E15") "  +2.19%  "
Set-TextValue $ws.Range("D16") "64.33"
Set-TextValue $ws.Range("E16") "  -1.78%  "
Set-TextValue $ws.Range("D17") "26.082.35"
Set-TextValue $ws.Range("E17") "  -2.00%  "
Set-TextValue $ws.Range("D18") "4.942"
Set-TextValue $ws.Range("E18") "  -1.36%  "
Set-TextValue $ws.Range("D19") "1.004"
Set-TextValue $ws.Range("E19") "  -0.12%  "
Set-TextValue $ws.Range("D20") "10.76"
Set-TextValue $ws.Range("E20") "  -2.05%  "
Set-TextValue $ws.Range("D21") "190.64"
Set-TextValue $ws.Range("E21") "  +2.55%  "
Set-TextValue $ws.Range("D22") "6.228"
Set-TextValue $ws.Range("E22") "  -0.62%  "
Set-TextValue $ws.Range("D23") "1.006"
Set-TextValue $ws.Range("E23") "  -0.01%  "
Set-TextValue $ws.Range("D24") "145.24"
Set-TextValue $ws.Range("E24") "  +0.49%  "
Set-TextValue $ws.Range("D25") "7.612"
Set-TextValue $ws.Range("E25") "  +1.45%  "
Set-TextValue $ws.Range("E26") "  +3.82%  "
Set-TextValue $ws.Range("E27") "  -0.21%  "
Set-TextValue $ws.Range("D28") "0.06675"
Set-TextValue $ws.Range("E28") "  +18.21%  "
Set-TextValue $ws.Range("D29") "1.316"
Set-TextValue $ws.Range("E29") "  -1.87%  "
Set-TextValue $ws.Range("D30") "1.315"
Set-TextValue $ws.Range("E30") "  -1.27%  "
Set-TextValue $ws.Range("D31") "3.537"
Set-TextValue $ws.Range("E31") "  +0.56%  "
Set-TextValue $ws.Range("D32") "3.517"
Set-TextValue $ws.Range("E32") "  +0.72%  "
Set-TextValue $ws.Range("D33") "1.649"
Set-TextValue $ws.Range("E33") "  +0.45%  "
Set-TextValue $ws.Range("D34") "1.017"
Set-TextValue $ws.Range("E34") "  -0.20%  "
Set-TextValue $ws.Range("D35") "0.6100"
Set-TextValue $ws.Range("E35") "  +1.18%  "
Set-TextValue $ws.Range("D36") "2.368"
Set-TextValue $ws.Range("E36") "  +0.20%  "
Set-TextValue $ws.Range("D37") "2.714"
Set-TextValue $ws.Range("D38") "6.223"
Set-TextValue $ws.Range("E38") "  +6.41%  "
Set-TextValue $ws.Range("D39") "0.01602"
Set-TextValue $ws.Range("E39") "  -0.57%  "
Set-TextValue $ws.Range("D40") "1.086.08"
Set-TextValue $ws.Range("D41") "0.8648"
Set-TextValue $ws.Range("E41") "  +0.64%  "
Set-TextValue $ws.Range("E42") "  +0.67%  "
Set-TextValue $ws.Range("D43") "100.70"
Set-TextValue $ws.Range("E43") "  +1.36%  "
Set-TextValue $ws.Range("D44") "1.816.82"
Set-TextValue $ws.Range("E44") "  -1.92%  "
Set-TextValue $ws.Range("D45") "0.00000000115"
Set-TextValue $ws.Range("E45") "  +4.45%  "
Set-TextValue $ws.Range("D46") "56.28"
Set-TextValue $ws.Range("E46") "  -0.73%  "
Set-TextValue $ws.Range("D47") "1.009"
Set-TextValue $ws.Range("E47") "  +0.46%  "
Set-TextValue $ws.Range("D48") "8.071"
Set-TextValue $ws.Range("E48") "  -0.89%  "
Set-TextValue $ws.Range("D49") "0.05241"
Set-TextValue $ws.Range("E49") "  +0.02%  "
Set-TextValue $ws.Range("E50") "  -0.95%  "
Set-TextValue $ws.Range("D51") "6.043"
Set-TextValue $ws.Range("E51") "  +4.59%  "
